$wb = $excel.ActiveWorkbook

$sheet = $wb.Worksheets.Item("Valve_30.0_600_3")
$sheet.Name = "Valve_33.0_600_3"

$valveList = $wb.Worksheets.Item("Valve List")
$valveList.Range("A10").Value = 33
